# Auto-generated edit script: update Exodus_Profits value columns per scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 927.9149
$ws.Range("I15").Value = 927.9149
$ws.Range("K15").Value = 2783.7447
$ws.Range("M15").Value = -2614.7447
$ws.Range("H17").Value = 2545.2727
$ws.Range("J17").Value = 2545.2727
$ws.Range("L17").Value = 7635.8181
$ws.Range("N17").Value = -7971.8181
$ws.Range("H58").Value = 160.3
$ws.Range("I58").Value = 121.125
$ws.Range("J58").Value = 317
$ws.Range("K58").Value = 363.375
$ws.Range("L58").Value = 951
$ws.Range("M58").Value = -213.375
$ws.Range("N58").Value = -1251
$ws.Range("H88").Value = 789.7143
$ws.Range("J88").Value = 821
$ws.Range("L88").Value = 821
$ws.Range("N88").Value = -1633
$ws.Range("H91").Value = 789.7143
$ws.Range("J91").Value = 821
$ws.Range("L91").Value = 821
$ws.Range("N91").Value = -3629
$ws.Range("H96").Value = 7419.2
$ws.Range("I96").Value = 7066
$ws.Range("K96").Value = 21198
$ws.Range("M96").Value = -19825
$ws.Range("H110").Value = 67993.336
$ws.Range("J110").Value = 67993.336
$ws.Range("L110").Value = 67993.336
$ws.Range("N110").Value = -76173.336
$ws.Range("H111").Value = 2029.4286
$ws.Range("J111").Value = 2055.4
$ws.Range("L111").Value = 6166.200000000001
$ws.Range("N111").Value = -12300.2
$ws.Range("H131").Value = 1029
$ws.Range("I131").Value = 1019.36365
$ws.Range("J131").Value = 1135
$ws.Range("K131").Value = 3058.09095
$ws.Range("L131").Value = 3405
$ws.Range("M131").Value = 1981.90905
$ws.Range("N131").Value = -13485
$ws.Range("H132").Value = 1174.963
$ws.Range("I132").Value = 1188.8679
$ws.Range("K132").Value = 3566.6037
$ws.Range("M132").Value = -1036.6037
$ws.Range("H136").Value = 60816
$ws.Range("J136").Value = 81632
$ws.Range("L136").Value = 81632
$ws.Range("N136").Value = -91832
$ws.Range("H138").Value = 2968.1428
$ws.Range("I138").Value = 1986.3704
$ws.Range("J138").Value = 3584.6047
$ws.Range("K138").Value = 5959.1112
$ws.Range("L138").Value = 10753.8141
$ws.Range("M138").Value = -819.1112000000003
$ws.Range("N138").Value = -21033.8141

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 13938.5
$ws.Range("I45").Value = 15701.286
$ws.Range("K45").Value = 15701.286
$ws.Range("M45").Value = -15324.286
$ws.Range("H61").Value = 2983.0356
$ws.Range("I61").Value = 2841
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 2841
$ws.Range("L61").Value = 4166.6665
$ws.Range("M61").Value = -2629
$ws.Range("N61").Value = -4590.6665
$ws.Range("H74").Value = 1964.375
$ws.Range("I74").Value = 1610
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 1610
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -736
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 1964.375
$ws.Range("I77").Value = 1610
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 8050
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -3682
$ws.Range("N77").Value = -26236
$ws.Range("H121").Value = 39993.4
$ws.Range("J121").Value = 39993.4
$ws.Range("L121").Value = 39993.4
$ws.Range("N121").Value = -43487.4
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 2267.4614
$ws.Range("I132").Value = 1972.8
$ws.Range("J132").Value = 3249.6667
$ws.Range("K132").Value = 5918.4
$ws.Range("L132").Value = 9749.000100000001
$ws.Range("M132").Value = -3388.4
$ws.Range("N132").Value = -14809.0001
$ws.Range("H136").Value = 2983.0356
$ws.Range("I136").Value = 2841
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 8523
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -5973
$ws.Range("N136").Value = -17599.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3408.3125
$ws.Range("I20").Value = 2193.2222
$ws.Range("J20").Value = 4970.5713
$ws.Range("K20").Value = 2193.2222
$ws.Range("L20").Value = 4970.5713
$ws.Range("M20").Value = -1946.2222
$ws.Range("N20").Value = -5464.5713
$ws.Range("H55").Value = 58547.25
$ws.Range("J55").Value = 65160
$ws.Range("L55").Value = 65160
$ws.Range("N55").Value = -65706
$ws.Range("H99").Value = 1954579
$ws.Range("I99").Value = 1431.8889
$ws.Range("J99").Value = 4465768
$ws.Range("K99").Value = 1431.8889
$ws.Range("L99").Value = 4465768
$ws.Range("M99").Value = 66.11110000000008
$ws.Range("N99").Value = -4468764
$ws.Range("H110").Value = 60991.668
$ws.Range("J110").Value = 60991.668
$ws.Range("L110").Value = 60991.668
$ws.Range("N110").Value = -69171.66800000001
$ws.Range("H114").Value = 92129
$ws.Range("J114").Value = 92213.664
$ws.Range("L114").Value = 92213.664
$ws.Range("N114").Value = -100891.664
$ws.Range("H119").Value = 83158.336
$ws.Range("J119").Value = 83158.336
$ws.Range("L119").Value = 83158.336
$ws.Range("N119").Value = -92834.336
$ws.Range("H132").Value = 32199.334
$ws.Range("J132").Value = 32199.334
$ws.Range("L132").Value = 32199.334
$ws.Range("N132").Value = -42319.334
$ws.Range("H134").Value = 1452.75
$ws.Range("I134").Value = 1275.1562
$ws.Range("K134").Value = 3825.4686
$ws.Range("M134").Value = -1290.4686
$ws.Range("H140").Value = 61902.637
$ws.Range("J140").Value = 61902.637
$ws.Range("L140").Value = 61902.637
$ws.Range("N140").Value = -72262.637

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3127849.8
$ws.Range("I99").Value = 2899.8
$ws.Range("J99").Value = 6252800
$ws.Range("K99").Value = 2899.8
$ws.Range("L99").Value = 6252800
$ws.Range("M99").Value = -1401.8
$ws.Range("N99").Value = -6255796
$ws.Range("H114").Value = 81713.5
$ws.Range("J114").Value = 81713.5
$ws.Range("L114").Value = 81713.5
$ws.Range("N114").Value = -90391.5
$ws.Range("H126").Value = 3127849.8
$ws.Range("I126").Value = 2899.8
$ws.Range("J126").Value = 6252800
$ws.Range("K126").Value = 8699.400000000001
$ws.Range("L126").Value = 18758400
$ws.Range("M126").Value = -6229.400000000001
$ws.Range("N126").Value = -18763340
$ws.Range("H134").Value = 1516.6938
$ws.Range("J134").Value = 2360.5715
$ws.Range("L134").Value = 7081.7145
$ws.Range("N134").Value = -12151.7145
$ws.Range("H141").Value = 163812.4
$ws.Range("J141").Value = 163812.4
$ws.Range("L141").Value = 163812.4
$ws.Range("N141").Value = -174172.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5299.9
$ws.Range("I3").Value = 1142.7142
$ws.Range("K3").Value = 3428.1426
$ws.Range("M3").Value = -3316.1426
$ws.Range("H122").Value = 1011838.9
$ws.Range("J122").Value = 1445014.2
$ws.Range("L122").Value = 13005127.8
$ws.Range("N122").Value = -13010027.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 38996.8
$ws.Range("J20").Value = 38996.8
$ws.Range("L20").Value = 38996.8
$ws.Range("N20").Value = -39486.8
$ws.Range("H96").Value = 18124.25
$ws.Range("J96").Value = 19086.666
$ws.Range("L96").Value = 19086.666
$ws.Range("N96").Value = -24578.666
$ws.Range("H108").Value = 65086.25
$ws.Range("J108").Value = 65241.43
$ws.Range("L108").Value = 65241.43
$ws.Range("N108").Value = -72921.42999999999
$ws.Range("H140").Value = 49150.77
$ws.Range("J140").Value = 49196
$ws.Range("L140").Value = 49196
$ws.Range("N140").Value = -59556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1258125
$ws.Range("I20").Value = 8000
$ws.Range("K20").Value = 8000
$ws.Range("M20").Value = -7774
$ws.Range("H136").Value = 7032.577
$ws.Range("I136").Value = 8062.9414
$ws.Range("J136").Value = 5086.3335
$ws.Range("K136").Value = 24188.8242
$ws.Range("L136").Value = 15259.0005
$ws.Range("M136").Value = -21638.8242
$ws.Range("N136").Value = -20359.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 100000000
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H122").Value = 3997.8823
$ws.Range("I122").Value = 2947.1365
$ws.Range("J122").Value = 5924.25
$ws.Range("K122").Value = 8841.4095
$ws.Range("L122").Value = 17772.75
$ws.Range("M122").Value = -6391.4095
$ws.Range("N122").Value = -22672.75
$ws.Range("H136").Value = 3435.8262
$ws.Range("I136").Value = 2966.2693
$ws.Range("J136").Value = 4046.25
$ws.Range("K136").Value = 8898.8079
$ws.Range("L136").Value = 12138.75
$ws.Range("M136").Value = -6348.8079
$ws.Range("N136").Value = -17238.75
$ws.Range("H137").Value = 20000
$ws.Range("J137").Value = 20000
$ws.Range("L137").Value = 20000
$ws.Range("N137").Value = -30200
